$d = $word.ActiveDocument

# Merge the stray split run in the "32-bit subtractor" cell so the
# trailing " (no saturation)" is a single run (matches "32-bit adder" row).
$d.Content.Find.Execute("32-bit subtractor (no saturation)", $true, $false, $false, $false, $false,
                         $true, 1, $false, "32-bit subtractor (no saturation)", 2)

$tbl = $d.Tables.Item(1)
Write-Output ("RowCount before: " + $tbl.Rows.Count)
$newRow = $tbl.Rows.Add()
Write-Output ("RowCount after: " + $tbl.Rows.Count)
$newRow.Cells.Item(1).Range.Text = "If-then (mux with register)"
$newRow.Cells.Item(2).Range.Text = "14"
$newRow.Range.LanguageID = 1033
